$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 233, pushing the existing data rows
# (old 233-245) down to 234-246.
$ws.Rows.Item(233).Insert()

# Populate the newly inserted row with the new weekly record. Most
# columns repeat the boilerplate from the surrounding rows (same market,
# region, product, etc.) while D/M/N/O/P/S carry the new week's figures.
$ws.Cells.Item(233, 1).Value = 10
$ws.Cells.Item(233, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(233, 3).Value = "La Araucanía"
$ws.Cells.Item(233, 4).Value = 44516
$ws.Cells.Item(233, 5).Value = 9
$ws.Cells.Item(233, 6).Value = "Fruta"
$ws.Cells.Item(233, 7).Value = 100108
$ws.Cells.Item(233, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(233, 9).Value = 100108002
$ws.Cells.Item(233, 10).Value = "Mango"
$ws.Cells.Item(233, 11).Value = "Sin especificar"
$ws.Cells.Item(233, 12).Value = "Primera"
$ws.Cells.Item(233, 13).Value = 155
$ws.Cells.Item(233, 14).Value = 8000
$ws.Cells.Item(233, 15).Value = 8000
$ws.Cells.Item(233, 16).Value = 8000
$ws.Cells.Item(233, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(233, 18).Value = "Perú"
$ws.Cells.Item(233, 19).Value = 2000
$ws.Cells.Item(233, 20).Value = 4
